$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 55901.5
$ws.Range("J57").Value = 55901.5
$ws.Range("L57").Value = 167704.5
$ws.Range("N57").Value = -168702.5
$ws.Range("H96").Value = 2478.182
$ws.Range("I96").Value = 806.8333
$ws.Range("K96").Value = 2420.4999
$ws.Range("M96").Value = -1047.4999
$ws.Range("H140").Value = 84550.42999999999
$ws.Range("J140").Value = 84475.836
$ws.Range("L140").Value = 84475.836
$ws.Range("N140").Value = -94835.836

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1764
$ws.Range("J2").Value = 2250
$ws.Range("L2").Value = 2250
$ws.Range("N2").Value = -2476
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H61").Value = 5709.8276
$ws.Range("I61").Value = 4813.091
$ws.Range("K61").Value = 4813.091
$ws.Range("M61").Value = -4601.091
$ws.Range("H95").Value = 54197.4
$ws.Range("J95").Value = 54197.4
$ws.Range("L95").Value = 54197.4
$ws.Range("N95").Value = -59689.4
$ws.Range("H97").Value = 2171.0454
$ws.Range("I97").Value = 1360.0588
$ws.Range("K97").Value = 1360.0588
$ws.Range("M97").Value = -864.0588
$ws.Range("H102").Value = 4297.2354
$ws.Range("I102").Value = 2727.75
$ws.Range("J102").Value = 8064
$ws.Range("K102").Value = 2727.75
$ws.Range("L102").Value = 8064
$ws.Range("M102").Value = -1105.75
$ws.Range("N102").Value = -11308
$ws.Range("H116").Value = 1764
$ws.Range("J116").Value = 2250
$ws.Range("L116").Value = 2250
$ws.Range("N116").Value = -6838
$ws.Range("H122").Value = 4587.1934
$ws.Range("I122").Value = 4598.2
$ws.Range("J122").Value = 4567.1816
$ws.Range("K122").Value = 13794.6
$ws.Range("L122").Value = 13701.5448
$ws.Range("M122").Value = -11344.6
$ws.Range("N122").Value = -18601.5448
$ws.Range("H132").Value = 8662.333000000001
$ws.Range("I132").Value = 4664.143
$ws.Range("K132").Value = 13992.429
$ws.Range("M132").Value = -11462.429
$ws.Range("H133").Value = 85990
$ws.Range("J133").Value = 85990
$ws.Range("L133").Value = 85990
$ws.Range("N133").Value = -91050
$ws.Range("H136").Value = 5709.8276
$ws.Range("I136").Value = 4813.091
$ws.Range("K136").Value = 14439.273
$ws.Range("M136").Value = -11889.273

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1764
$ws.Range("J3").Value = 2250
$ws.Range("L3").Value = 2250
$ws.Range("N3").Value = -2478
$ws.Range("H20").Value = 1464.0312
$ws.Range("I20").Value = 1295.1052
$ws.Range("J20").Value = 1710.9231
$ws.Range("K20").Value = 1295.1052
$ws.Range("L20").Value = 1710.9231
$ws.Range("M20").Value = -1048.1052
$ws.Range("N20").Value = -2204.9231
$ws.Range("H22").Value = 482.81818
$ws.Range("I22").Value = 438.875
$ws.Range("K22").Value = 438.875
$ws.Range("M22").Value = -265.875
$ws.Range("H107").Value = 7432.04
$ws.Range("I107").Value = 7202.8335
$ws.Range("K107").Value = 7202.8335
$ws.Range("M107").Value = -5282.8335
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 242966.33
$ws.Range("J9").Value = 242966.33
$ws.Range("L9").Value = 242966.33
$ws.Range("N9").Value = -243302.33
$ws.Range("H32").Value = 4185.875
$ws.Range("I32").Value = 3355.5715
$ws.Range("K32").Value = 3355.5715
$ws.Range("M32").Value = -3039.5715
$ws.Range("H58").Value = 2009.579
$ws.Range("J58").Value = 3950
$ws.Range("L58").Value = 3950
$ws.Range("N58").Value = -4356
$ws.Range("H136").Value = 2009.579
$ws.Range("J136").Value = 3950
$ws.Range("L136").Value = 11850
$ws.Range("N136").Value = -16950

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 460.83334
$ws.Range("I13").Value = 479.66666
$ws.Range("J13").Value = 442
$ws.Range("K13").Value = 1438.99998
$ws.Range("L13").Value = 1326
$ws.Range("M13").Value = -1270.99998
$ws.Range("N13").Value = -1662
$ws.Range("H92").Value = 207.36363
$ws.Range("I92").Value = 171.8
$ws.Range("J92").Value = 237
$ws.Range("K92").Value = 515.4000000000001
$ws.Range("L92").Value = 711
$ws.Range("M92").Value = 732.5999999999999
$ws.Range("N92").Value = -3207
$ws.Range("H107").Value = 5425.3335
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 5425.3335
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 16276.0005
$ws.Range("N107").Value = -20116.0005
$ws.Range("H113").Value = 756.82355
$ws.Range("I113").Value = 416.8
$ws.Range("J113").Value = 898.5
$ws.Range("K113").Value = 1250.4
$ws.Range("L113").Value = 2695.5
$ws.Range("M113").Value = 919.5999999999999
$ws.Range("N113").Value = -7035.5
$ws.Range("H131").Value = 1730.5366
$ws.Range("I131").Value = 1275
$ws.Range("J131").Value = 1808.6285
$ws.Range("K131").Value = 3825
$ws.Range("L131").Value = 5425.8855
$ws.Range("M131").Value = 1215
$ws.Range("N131").Value = -15505.8855
$ws.Range("H132").Value = 2158.8572
$ws.Range("I132").Value = 1498
$ws.Range("J132").Value = 2238.16
$ws.Range("K132").Value = 13482
$ws.Range("L132").Value = 20143.44
$ws.Range("M132").Value = -10952
$ws.Range("N132").Value = -25203.44

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H113").Value = 2957.4443
$ws.Range("I113").Value = 3753.7144
$ws.Range("K113").Value = 3753.7144
$ws.Range("M113").Value = -1583.7144
$ws.Range("H132").Value = 7304.7026
$ws.Range("J132").Value = 4783.3335
$ws.Range("L132").Value = 14350.0005
$ws.Range("N132").Value = -19410.0005

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3584.0625
$ws.Range("I22").Value = 501
$ws.Range("J22").Value = 3789.6
$ws.Range("K22").Value = 501
$ws.Range("L22").Value = 3789.6
$ws.Range("M22").Value = -206
$ws.Range("N22").Value = -4379.6
$ws.Range("H27").Value = 3584.0625
$ws.Range("I27").Value = 501
$ws.Range("J27").Value = 3789.6
$ws.Range("K27").Value = 501
$ws.Range("L27").Value = 3789.6
$ws.Range("M27").Value = -394
$ws.Range("N27").Value = -4003.6
$ws.Range("H61").Value = 23966.137
$ws.Range("I61").Value = 18721.354
$ws.Range("K61").Value = 18721.354
$ws.Range("M61").Value = -18519.354
$ws.Range("H92").Value = 82096
$ws.Range("J92").Value = 82096
$ws.Range("L92").Value = 82096
$ws.Range("N92").Value = -87088
$ws.Range("H100").Value = 3790
$ws.Range("I100").Value = 3485.7144
$ws.Range("K100").Value = 3485.7144
$ws.Range("M100").Value = -2944.7144
$ws.Range("H113").Value = 23966.137
$ws.Range("I113").Value = 18721.354
$ws.Range("K113").Value = 18721.354
$ws.Range("M113").Value = -16551.354
$ws.Range("H136").Value = 3825.5557
$ws.Range("I136").Value = 3853.7812
$ws.Range("K136").Value = 11561.3436
$ws.Range("M136").Value = -9011.3436

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 75000
$ws.Range("J75").Value = 75000
$ws.Range("L75").Value = 75000
$ws.Range("N75").Value = -76872
$ws.Range("H78").Value = 75000
$ws.Range("J78").Value = 75000
$ws.Range("L78").Value = 225000
$ws.Range("N78").Value = -234360
$ws.Range("H100").Value = 3766.111
$ws.Range("I100").Value = 1547.5
$ws.Range("K100").Value = 3095
$ws.Range("M100").Value = -2554
$ws.Range("H113").Value = 3473135.5
$ws.Range("I113").Value = 6944867
$ws.Range("K113").Value = 20834601
$ws.Range("M113").Value = -20832431
$ws.Range("H114").Value = 51500
$ws.Range("J114").Value = 51500
$ws.Range("L114").Value = 51500
$ws.Range("N114").Value = -60178
$ws.Range("H136").Value = 4072.1892
$ws.Range("I136").Value = 2731.4092
$ws.Range("K136").Value = 8194.2276
$ws.Range("M136").Value = -5644.2276
